$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("FID", "Name", "Age", "Phone", "Vaccine_Dose")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match column widths observed in the target workbook
# (iron_native stores width = ColumnWidth + 5/6, so back out the COM value)
$ws.Columns.Item(1).ColumnWidth = 9.893229166666666
$ws.Columns.Item(2).ColumnWidth = 11.893229166666666
$ws.Columns.Item(3).ColumnWidth = 10.436197916666666
$ws.Columns.Item(4).ColumnWidth = 13.529947916666666
$ws.Columns.Item(5).ColumnWidth = 13.436197916666666

$ws.Range("A2").Select()
